$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
  2 = @{ "J" = 3.2; "M" = 1.12; "N" = 2.46; "Q" = 2.66; "T" = 2.38; "AC" = 8.800000000000001; "AD" = 32 }
  3 = @{ "G" = 3.45; "J" = 3.55; "K" = 3.95; "L" = 1.3; "N" = 4.1; "Q" = 1.79; "R" = 1.41; "S" = 2.94; "T" = 1.67; "U" = 2.26; "W" = 1.4; "X" = 21; "Y" = 12.5; "AB" = 15.5; "AE" = 29; "AG" = 17; "AI" = 42; "AK" = 1000; "AL" = 50; "AN" = 30; "AO" = 17 }
  4 = @{ "F" = 2.74; "H" = 2.5; "I" = 3.1; "J" = 3; "L" = 1.46; "Q" = 2; "S" = 2.16; "V" = 1.48; "W" = 1.4; "Y" = 1000; "AB" = 1000 }
  5 = @{ "F" = 1.64; "G" = 1.79; "H" = 7.2; "J" = 3.2; "K" = 3.6; "L" = 1.49; "N" = 2.44; "O" = 1.57; "Q" = 2.7; "S" = 5.8; "T" = 2.48; "U" = 1.54; "W" = 2.28; "X" = 9.800000000000001; "Y" = 20; "AF" = 10 }
  6 = @{ "F" = 4.5; "G" = 4.9; "H" = 1.87; "I" = 1.91; "J" = 3.75; "L" = 1.44; "N" = 3.4; "O" = 1.38; "P" = 1.81; "Q" = 2.14; "T" = 1.98; "U" = 1.92; "V" = 2.08; "W" = 1.26; "Z" = 10.5; "AB" = 15; "AF" = 34; "AJ" = 120; "AN" = 85; "AO" = 15.5 }
  7 = @{ "F" = 1.28; "G" = 1.35; "H" = 16; "I" = 22; "K" = 6.2; "M" = 1.07; "N" = 3.2; "O" = 1.35; "P" = 1.74; "Q" = 2.02; "R" = 1.28; "T" = 2.74; "U" = 1.44; "V" = 1.05; "W" = 3.85; "X" = 14; "Y" = 1000; "AB" = 6.4; "AC" = 14; "AF" = 6.4; "AG" = 13; "AH" = 55; "AJ" = 9.6; "AK" = 21; "AN" = 7.8 }
  8 = @{ "F" = 1.51; "I" = 8.199999999999999; "J" = 3.8; "Q" = 1.93 }
  9 = @{ "F" = 1.54; "G" = 1.55; "H" = 6.8; "N" = 4.9; "O" = 1.23; "P" = 2.36; "S" = 2.72; "T" = 1.81; "W" = 2.8; "Y" = 980; "AA" = 190; "AE" = 90; "AH" = 21; "AM" = 100; "AO" = 85 }
  10 = @{ "F" = 1.74; "G" = 1.94; "H" = 4.7; "I" = 5.9; "J" = 3.5; "K" = 4.3; "L" = 1.38; "M" = 1.06; "N" = 3.45; "O" = 1.3; "P" = 1.84; "Q" = 1.9; "R" = 1.32; "S" = 3.3; "T" = 1.81; "U" = 1.94; "V" = 1.21; "W" = 2.06; "X" = 15.5; "Y" = 18.5; "AB" = 10; "AC" = 9.199999999999999; "AD" = 22; "AF" = 12; "AG" = 12.5; "AH" = 22; "AJ" = 22; "AK" = 21; "AN" = 14 }
  11 = @{ "F" = 2.5; "G" = 2.58; "H" = 3.45; "K" = 3.1; "L" = 1.59; "M" = 1.13; "N" = 2.74; "O" = 1.55; "P" = 1.55; "Q" = 2.66; "S" = 5.4; "T" = 2.1; "V" = 1.38; "W" = 1.64; "Z" = 21; "AE" = 310; "AL" = 70; "AM" = 190 }
  12 = @{ "F" = 1.95; "G" = 2.16; "H" = 4.2; "I" = 5.6; "P" = 1.57 }
  13 = @{ "F" = 2.34; "G" = 2.38; "H" = 3.55; "I" = 3.7; "J" = 3.3; "K" = 3.45; "L" = 1.56; "M" = 1.13; "N" = 2.66; "O" = 1.56; "P" = 1.56; "Q" = 2.66; "R" = 1.19; "S" = 5.6; "T" = 2.22; "U" = 1.71; "V" = 1.37; "W" = 1.72; "X" = 8.800000000000001; "Y" = 10; "Z" = 23; "AA" = 80; "AB" = 7.4; "AC" = 8; "AD" = 16.5; "AE" = 60; "AF" = 12; "AG" = 12; "AH" = 30; "AI" = 85; "AJ" = 34; "AK" = 34; "AL" = 65; "AM" = 190; "AN" = 34; "AO" = 80 }
  14 = @{ "F" = 1.63; "G" = 1.69; "H" = 5.6; "I" = 7; "J" = 3.9; "K" = 4.3; "L" = 1.33; "N" = 3.6; "O" = 1.32; "P" = 1.94; "Q" = 1.91; "R" = 1.36; "T" = 1.81; "U" = 1.9; "V" = 1.18; "W" = 2.44; "AB" = 9.4; "AF" = 11.5; "AJ" = 19.5; "AN" = 12.5; "AO" = 160 }
  15 = @{ "H" = 3.55; "I" = 3.8; "J" = 3.8; "Q" = 1.8; "S" = 3; "U" = 2.22; "V" = 1.36; "Y" = 16.5; "AC" = 9.199999999999999 }
}

foreach ($r in $changes.Keys) {
    foreach ($c in $changes[$r].Keys) {
        $addr = "$c$r"
        $ws.Range($addr).Value = $changes[$r][$c]
    }
}
